# Fixed update to excel issue
# Apply corrected values in column C ("My 4 Weeks Forecast") for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    9  = 456
    11 = 780.8000000000001
    14 = 2408
    15 = 166.5
    21 = 467
    26 = 1065.9
    30 = 57.00000000000001
    34 = 12.78
    38 = 40
    49 = 58.16000000000001
    52 = 216
    76 = 115
}

foreach ($row in $updates.Keys) {
    $ws.Range("C$row").Value = $updates[$row]
}
